$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Regime Atual
$ws.Range("B2").Value = 0.6200835989683386
$ws.Range("C2").Value = 0.1393566213130238
$ws.Range("D2").Value = 0.5510331394514802
$ws.Range("E2").Value = 0.247079201554608
$ws.Range("F2").Value = 299.1466261603687

# Row 3 - Nova Proposta
$ws.Range("B3").Value = 0.617708066649426
$ws.Range("C3").Value = 0.140169244854642
$ws.Range("D3").Value = 0.548308748648395
$ws.Range("E3").Value = 0.2417126824104117
$ws.Range("F3").Value = 337.3358639310761
$ws.Range("G3").Value = 38.18923777070739

# Row 4 - Nova c/ Aliq. Máxima
$ws.Range("B4").Value = 0.6170084104030436
$ws.Range("C4").Value = 0.1404290002923165
$ws.Range("D4").Value = 0.5474716944229102
$ws.Range("E4").Value = 0.2403055863644433
$ws.Range("F4").Value = 349.5950229090772
$ws.Range("G4").Value = 50.44839674870849
